# Generate Report for Handoff
# Refresh the localization-status report: the handoff run for rows 4-7
# (the 2d015710 / a5751771 / bc4491a8 / e6a30c27 files) has finished, so
# their Priority flips from "low" to "ht" and the Latest Handoff
# Datetime is bumped to the new run's timestamp, for both the zh-cn and
# de-de target-language sheets. The Overview sheet's "Latest HO Xliff
# Generate Date" column mirrors the de-de handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-17 08:30:48"

    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-17 08:30:43"

    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-17 08:30:48"
}
